# Bugfixed the naive forecaster component module
# The first forecast row (row 2, date 2007-11-14 / y_0=2007) was dropped, shifting
# every remaining row up by one, and the y_0_forecast / y_1_forecast columns were
# recomputed for the AR(2) model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete first data row; this shifts all following rows up by one
# and updates the sheet dimension automatically (A1:E19 -> A1:E18).
$ws.Rows("2").Delete()

# Correct a couple of slightly-changed y_0_forecast (column C) values.
$ws.Range("C3").Value = -1.324983933426893
$ws.Range("C5").Value = -0.29958481534893

# The y_1_forecast column (E) no longer has values for the first few rows.
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E5").ClearContents()

# Recomputed y_1_forecast values for the remaining rows.
$ws.Range("E6").Value = -0.4223781730902543
$ws.Range("E7").Value = -0.3577371449824729
$ws.Range("E8").Value = -0.2877346565283379
$ws.Range("E9").Value = -0.2097319935285391
$ws.Range("E10").Value = -0.3036570471216304
$ws.Range("E11").Value = -0.1907914728172644
$ws.Range("E12").Value = -0.2529765062333933
$ws.Range("E13").Value = -0.428077259747528
$ws.Range("E14").Value = -0.2112001730687485
$ws.Range("E15").Value = -1.701252732314051
$ws.Range("E16").Value = -0.6241481568271312
$ws.Range("E17").Value = -0.2119687890143274
$ws.Range("E18").Value = -0.06564014165270082
